$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "24.071.43"
$ws.Range("E2").Value = "  -3.25%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.638.71"
$ws.Range("E3").Value = "  -2.92%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  -0.55%  "

# Row 5 - USDC
Set-TextValue "D5" "1.001"
$ws.Range("E5").Value = "  -0.50%  "

# Row 6 - BNB
Set-TextValue "D6" "306.67"
$ws.Range("E6").Value = "  -2.82%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.74%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3822"
$ws.Range("E8").Value = "  -3.86%  "

# Row 9 - BinanceUSD
Set-TextValue "D9" "1.001"
$ws.Range("E9").Value = "  -0.63%  "

# Row 10 - OKB
Set-TextValue "D10" "49.14"
$ws.Range("E10").Value = "  -6.24%  "

# Row 11 - Polygon
Set-TextValue "D11" "1.340"
$ws.Range("E11").Value = "  -6.73%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.08383"
$ws.Range("E12").Value = "  -3.85%  "

# Row 13 - Solana
Set-TextValue "D13" "23.59"
$ws.Range("E13").Value = "  -7.33%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.051"
$ws.Range("E14").Value = "  -4.33%  "

# Row 15 - ShibaInu
Set-TextValue "D15" "0.00001272"
$ws.Range("E15").Value = "  -4.70%  "

# Row 16 - Chainlink
Set-TextValue "D16" "7.414"
$ws.Range("E16").Value = "  -5.61%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "1.642.64"
$ws.Range("E17").Value = "  +7.59%  "

# Row 18 - Litecoin
Set-TextValue "D18" "95.29"
$ws.Range("E18").Value = "  +0.71%  "

# Row 19 - TRON
Set-TextValue "D19" "0.06868"
$ws.Range("E19").Value = "  -5.23%  "

# Row 20 - Avalanche
Set-TextValue "D20" "20.29"
$ws.Range("E20").Value = "  -0.52%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.863"
$ws.Range("E21").Value = "  -4.16%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.58%  "

# Row 23 - Cosmos
Set-TextValue "D23" "13.47"

# Row 24 - WrappedBTC
Set-TextValue "D24" "24.080.12"
$ws.Range("E24").Value = "  -3.17%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.328"
$ws.Range("E25").Value = "  -3.01%  "

# Row 26 - LidoDAOToken
Set-TextValue "D26" "2.683"
$ws.Range("E26").Value = "  -4.93%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "22.26"
$ws.Range("E27").Value = "  -3.49%  "

# Row 28 - now Filecoin (was Monero)
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D28" "8.750"
$ws.Range("E28").Value = "  +8.51%  "

# Row 29 - now Monero (was Filecoin)
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D29" "156.78"
$ws.Range("E29").Value = "  -2.85%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "139.72"
$ws.Range("E30").Value = "  -5.69%  "

# Row 31 - HuobiToken
Set-TextValue "D31" "5.324"
$ws.Range("E31").Value = "  -11.41%  "

# Row 32 - WEMIXTOKEN
Set-TextValue "D32" "2.420"
$ws.Range("E32").Value = "  -7.10%  "

# Row 33 - WrappedliquidstakedEther2.0
Set-TextValue "D33" "1.822.38"
$ws.Range("E33").Value = "  -17.18%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "6.877"
$ws.Range("E34").Value = "  -2.39%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.07964"
$ws.Range("E35").Value = "  -6.33%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.02872"
$ws.Range("E36").Value = "  -7.43%  "

# Row 37 - Algorand
Set-TextValue "D37" "0.2666"
$ws.Range("E37").Value = "  -6.63%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.9449"
$ws.Range("E38").Value = "  -8.08%  "

# Row 39 - Stellar
Set-TextValue "D39" "0.09151"
$ws.Range("E39").Value = "  -4.88%  "

# Row 40 - TrustWalletToken
Set-TextValue "D40" "1.449"
$ws.Range("E40").Value = "  -1.44%  "

# Row 41 - FraxShare
Set-TextValue "D41" "9.829"
$ws.Range("E41").Value = "  -9.25%  "

# Row 42 - TheSandbox
Set-TextValue "D42" "0.7487"
$ws.Range("E42").Value = "  -7.21%  "

# Row 43 - Aptos
Set-TextValue "D43" "12.96"
$ws.Range("E43").Value = "  -6.78%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "15.99"
$ws.Range("E44").Value = "  -5.82%  "

# Row 45 - Decentraland
Set-TextValue "D45" "0.6842"
$ws.Range("E45").Value = "  -5.77%  "

# Row 46 - NEARProtocol
Set-TextValue "D46" "2.448"
$ws.Range("E46").Value = "  -6.69%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "4.077"
$ws.Range("E47").Value = "  -3.35%  "

# Row 48 - Frax
Set-TextValue "D48" "1.001"
$ws.Range("E48").Value = "  -1.45%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.08332"
$ws.Range("E49").Value = "  -6.59%  "

# Row 50 - Flow
Set-TextValue "D50" "1.249"
$ws.Range("E50").Value = "  -9.20%  "

# Row 51 - Quant
Set-TextValue "D51" "132.00"
$ws.Range("E51").Value = "  -5.65%  "
